# Add data for 2024-12-24
# Updates K-column (year 2024) totals across "Citywide Totals", "By Neighborhood",
# and per-neighborhood sheets to reflect newly added crime records for 2024-12-24.
# A couple of I-column (year 2022) cells are also corrected where the source
# diff indicates a value change.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7644
$ws.Range("K3").Value = 7902
$ws.Range("I4").Value = 1814
$ws.Range("K4").Value = 1664
$ws.Range("K5").Value = 567
$ws.Range("K6").Value = 8803
$ws.Range("I7").Value = 26273
$ws.Range("K7").Value = 26580

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 194
$ws.Range("K7").Value = 796
$ws.Range("K8").Value = 1736
$ws.Range("K14").Value = 127
$ws.Range("K15").Value = 274
$ws.Range("K18").Value = 179
$ws.Range("K19").Value = 767
$ws.Range("K20").Value = 652
$ws.Range("K25").Value = 123
$ws.Range("K29").Value = 1465
$ws.Range("K31").Value = 314
$ws.Range("K33").Value = 1118
$ws.Range("K34").Value = 153
$ws.Range("K36").Value = 342
$ws.Range("K37").Value = 877
$ws.Range("K43").Value = 222
$ws.Range("K48").Value = 331
$ws.Range("K49").Value = 150
$ws.Range("K52").Value = 685
$ws.Range("K53").Value = 335
$ws.Range("K55").Value = 289
$ws.Range("K60").Value = 158
$ws.Range("I63").Value = 235
$ws.Range("K63").Value = 72
$ws.Range("K65").Value = 622
$ws.Range("K66").Value = 78
$ws.Range("K67").Value = 1032
$ws.Range("K69").Value = 62
$ws.Range("K76").Value = 367
$ws.Range("K78").Value = 329
$ws.Range("K79").Value = 652
$ws.Range("K83").Value = 562
$ws.Range("K84").Value = 216
$ws.Range("K85").Value = 1220
$ws.Range("K86").Value = 162
$ws.Range("K88").Value = 284
$ws.Range("K89").Value = 395
$ws.Range("K90").Value = 256
$ws.Range("K94").Value = 359
$ws.Range("K95").Value = 443
$ws.Range("K96").Value = 284
$ws.Range("K97").Value = 218
$ws.Range("K99").Value = 446
$ws.Range("I101").Value = 26273
$ws.Range("K101").Value = 26580

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K3").Value = 32
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 252
$ws.Range("K7").Value = 796

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 395

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 401
$ws.Range("K6").Value = 300
$ws.Range("K7").Value = 1220

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 190
$ws.Range("K7").Value = 685

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K4").Value = 18
$ws.Range("K6").Value = 138
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K6").Value = 582
$ws.Range("K7").Value = 1736

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 196
$ws.Range("K7").Value = 562

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 279
$ws.Range("K3").Value = 394
$ws.Range("K4").Value = 55
$ws.Range("K6").Value = 358
$ws.Range("K7").Value = 1118

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 149
$ws.Range("K6").Value = 103
$ws.Range("K7").Value = 443

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 249
$ws.Range("K3").Value = 291
$ws.Range("K6").Value = 265
$ws.Range("K7").Value = 877

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 202
$ws.Range("K4").Value = 28
$ws.Range("K6").Value = 231
$ws.Range("K7").Value = 622

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 189
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 446

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 94
$ws.Range("K6").Value = 122
$ws.Range("K7").Value = 314

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 372
$ws.Range("K7").Value = 1032

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 119
$ws.Range("K6").Value = 284

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K4").Value = 67
$ws.Range("K6").Value = 434
$ws.Range("K7").Value = 1465

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 223
$ws.Range("K6").Value = 259
$ws.Range("K7").Value = 767

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 182
$ws.Range("K7").Value = 367

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 100
$ws.Range("K3").Value = 84
$ws.Range("K7").Value = 329

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 84
$ws.Range("K7").Value = 289

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 214
$ws.Range("K3").Value = 205
$ws.Range("K7").Value = 652

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 223
$ws.Range("K3").Value = 203
$ws.Range("K7").Value = 652

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 128
$ws.Range("K3").Value = 108
$ws.Range("K7").Value = 342

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K2").Value = 60
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 167
$ws.Range("K7").Value = 359

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 274

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K2").Value = 43
$ws.Range("K7").Value = 218

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 89
$ws.Range("K5").Value = 5
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 69
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 158

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 65
$ws.Range("K7").Value = 222
